$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "35.111.07"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.893.79"
$ws.Range("E3").Value = "  +1.51%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.73"
$ws.Range("E5").Value = "  +2.53%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +5.85%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.21%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.19"
$ws.Range("E8").Value = "  -4.04%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.346"
$ws.Range("E9").Value = "  +4.60%  "

# Row 10: OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.74"
$ws.Range("E10").Value = "  +12.60%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0715"
$ws.Range("E11").Value = "  +2.58%  "

# Row 12: TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0992"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("E13").Value = "  +1.48%  "

# Row 14: Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.20"
$ws.Range("E14").Value = "  +5.39%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.693"
$ws.Range("E15").Value = "  +1.75%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "1.880.68"
$ws.Range("E16").Value = "  +1.93%  "

# Row 17: Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.78"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "35.098.42"
$ws.Range("E18").Value = "  -0.64%  "

# Row 19: Litecoin
$ws.Range("E19").Value = "  +2.74%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("E20").Value = "  +1.94%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "239.77"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22: Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.39"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23: Uniswap
$ws.Range("E23").Value = "  +1.01%  "

# Row 24: Dai
$ws.Range("E24").Value = "  -0.12%  "

# Row 25: PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +23.65%  "

# Row 26: Toncoin
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +0.88%  "

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.39"
$ws.Range("E27").Value = "  +0.42%  "

# Row 28: Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.44"
$ws.Range("E28").Value = "  +2.71%  "

# Row 29: EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.28"
$ws.Range("E29").Value = "  +2.78%  "

# Row 31: EURNeutrino
$ws.Range("D31").Value = "4.130.88"
$ws.Range("E31").Value = "  +20.99%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.11"
$ws.Range("E32").Value = "  +2.05%  "

# Row 33: Hedera
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0560"
$ws.Range("E33").Value = "  -0.74%  "

# Row 34: ImmutableX
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.938"
$ws.Range("E34").Value = "  +14.87%  "

# Row 35: BinanceUSD
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36: InternetComputer(DFINITY)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.06"
$ws.Range("E36").Value = "  +0.10%  "

# Row 37: WEMIXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.75"
$ws.Range("E37").Value = "  -5.32%  "

# Row 38: LidoDAOToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.02"
$ws.Range("E38").Value = "  -2.63%  "

# Row 39: TrustWalletToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  +0.98%  "

# Row 40: ARBITRUM
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.10"
$ws.Range("E40").Value = "  -1.52%  "

# Row 41: VeChain
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0207"
$ws.Range("E41").Value = "  +1.80%  "

# Row 42: InjectiveProtocol
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.04"
$ws.Range("E42").Value = "  +5.27%  "

# Row 43: Kaspa
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0629"
$ws.Range("E43").Value = "  +4.39%  "

# Row 44: Aave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "89.38"
$ws.Range("E44").Value = "  -1.86%  "

# Row 45: Maker
$ws.Range("D45").Value = "1.332.83"
$ws.Range("E45").Value = "  -1.19%  "

# Row 46: MultiversX
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.74"
$ws.Range("E46").Value = "  +39.50%  "

# Row 47: RenderToken
$ws.Range("E47").Value = "  +0.32%  "

# Row 48: HuobiToken
$ws.Range("E48").Value = "  -0.40%  "

# Row 49: MXToken
$ws.Range("E49").Value = "  +0.99%  "

# Row 50: FraxShare
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.45"
$ws.Range("E50").Value = "  -2.86%  "

# Row 51: RocketPoolETH
$ws.Range("D51").Value = "2.076.11"
$ws.Range("E51").Value = "  +1.26%  "
